# CIERRE 31 ENE 22
# Updates the "Hoja1" weekly payroll sheet with the closing figures for the
# week and moves the on-screen view/selection down to the last block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Payroll figures -------------------------------------------------
# Block 1 (rows 2-7): INFONAVIT discount column (K4) now has a value; the
# SUM(K3:K5) total in K6 recalculates automatically.
$ws.Range("K4").Value2 = 867

# Block 2 (rows 20-26): updated total in K21; SUM(K21:K23) in K24 follows.
$ws.Range("K21").Value2 = 2380

# Block 3 (rows 36-41): days worked (D38) and the corresponding pay (E38)
# are adjusted; SUM(E38:E40) in E41 follows automatically. The 3 extra
# hours previously recorded in J39 are removed (value cleared, format kept).
$ws.Range("D38").Value2 = 5
$ws.Range("E38").Value2 = 1833
$ws.Range("J39").ClearContents()

# --- Date stamps -------------------------------------------------------
# C14 drives I14/C32/I32/C48/I48/C65 via formulas already in the sheet;
# simply re-asserting the TODAY()-based formula chain is enough to refresh
# the cached date stamps for the closing.
$ws.Range("C14").Formula = "=TODAY()"
$ws.Range("I14").Formula = "=C14"
$ws.Range("C32").Formula = "=I14"
$ws.Range("I32").Formula = "=C32"
$ws.Range("C48").Formula = "=C32"
$ws.Range("I48").Formula = "=C48"
$ws.Range("C65").Formula = "=I48"

# --- View state ----------------------------------------------------------
# Scroll the window so row 40 is at the top and select H60, matching where
# the person closing the week left off.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H60").Select() | Out-Null
